# Update countries & provincias Spain
# - Insert "El Salvador" into its sorted position (row 92), pushing
#   Lituania / Nueva Zelanda / Somalia / Gabon down one row each, with
#   refreshed case numbers.
# - Refresh a couple of unrelated country rows (Alemania row 11,
#   Republica de Chipre row 110).
# - Bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Alemania (row 11): casos activos / recuperados update ---
$ws.Cells.Item(11, 4).Value = 156900
$ws.Cells.Item(11, 5).Value = 12734

# --- El Salvador enters the table at row 92, shifting the following
#     countries down by one row (values only; column A text moves too) ---
$ws.Cells.Item(92, 1).Value = "El Salvador"
$ws.Cells.Item(92, 2).Value = 1571
$ws.Cells.Item(92, 3).Value = 73
$ws.Cells.Item(92, 4).Value = 531
$ws.Cells.Item(92, 5).Value = 1009
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 1
$ws.Cells.Item(92, 8).Value = 31

$ws.Cells.Item(93, 1).Value = "Lituania"
$ws.Cells.Item(93, 2).Value = 1562
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 1025
$ws.Cells.Item(93, 5).Value = 477
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 60

$ws.Cells.Item(94, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(94, 2).Value = 1503
$ws.Cells.Item(94, 3).Value = 0
$ws.Cells.Item(94, 4).Value = 1447
$ws.Cells.Item(94, 5).Value = 35
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 21

$ws.Cells.Item(95, 1).Value = "Somalia"
$ws.Cells.Item(95, 2).Value = 1502
$ws.Cells.Item(95, 3).Value = 0
$ws.Cells.Item(95, 4).Value = 178
$ws.Cells.Item(95, 5).Value = 1265
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 59

$ws.Cells.Item(96, 1).Value = "Gabon"
$ws.Cells.Item(96, 2).Value = 1502
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 318
$ws.Cells.Item(96, 5).Value = 1172
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 12

# --- Republica de Chipre (row 110): casos activos / recuperados update ---
$ws.Cells.Item(110, 4).Value = 516
$ws.Cells.Item(110, 5).Value = 385

# --- Timestamp footer update (A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Mayo de 2020 a las 08:35"
